# Rename the embedded logo pictures that live in the document's headers
# and footers. This mirrors a straightforward "rename picture" edit done
# in Word's UI (right-click a picture -> Size & Properties / Alt Text,
# or the Selection Pane "Rename") which updates the drawing's display
# name (wp:docPr/@name, mirrored onto pic:cNvPr/@name by Word itself):
#
#   * BTec_Logo-Orange picture (in both headers)  : image1.jpg -> image2.jpg
#   * PearsonLogo picture      (in both footers)   : image2.png -> image1.png
#
# InlineShape objects don't reliably expose a usable .Name setter when
# addressed directly off a HeaderFooter.Range in this host, so each
# picture is selected first and then renamed through $word.Selection -
# exactly like a user would click the picture and rename it.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-HeaderFooterPicture($range, $newName) {
    $shape = $range.InlineShapes.Item(1)
    $shape.Select()
    $word.Selection.InlineShapes.Item(1).Name = $newName
}

# Headers: BTec_Logo-Orange -> image2.jpg
for ($h = 1; $h -le 2; $h++) {
    $hdr = $sec.Headers.Item($h)
    if ($hdr.Exists) {
        Rename-HeaderFooterPicture $hdr.Range "image2.jpg"
    }
}

# Footers: PearsonLogo -> image1.png
for ($f = 1; $f -le 2; $f++) {
    $ftr = $sec.Footers.Item($f)
    if ($ftr.Exists) {
        Rename-HeaderFooterPicture $ftr.Range "image1.png"
    }
}

Write-Output "Renamed header/footer logo pictures."
